# Updates crypto price/volume/hour data for rows 2-51 (Sheet1)
# per the "Updated symbol list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new literal text (kept as Text so values
# such as "10" or "-2.78%" are stored verbatim, matching the source data).
$updates = [ordered]@{
    'D2' = '326.52'
    'E2' = '-2.78%'
    'G2' = '10'
    'D3' = '44.56'
    'E3' = '1.71%'
    'G3' = '10'
    'D4' = '5.610'
    'E4' = '-2.57%'
    'G4' = '10'
    'D5' = '0.08057'
    'E5' = '-3.12%'
    'G5' = '10'
    'D6' = '1.905'
    'E6' = '-2.91%'
    'G6' = '10'
    'D7' = '4.304'
    'E7' = '-4.90%'
    'G7' = '10'
    'D8' = '2.721'
    'E8' = '-5.82%'
    'G8' = '10'
    'D9' = '0.9446'
    'G9' = '10'
    'E10' = '-6.39%'
    'G10' = '10'
    'D11' = '0.1848'
    'E11' = '-6.12%'
    'G11' = '10'
    'D12' = '0.09868'
    'E12' = '-7.28%'
    'G12' = '10'
    'D13' = '0.04252'
    'E13' = '-7.87%'
    'G13' = '10'
    'D14' = '0.1066'
    'E14' = '-0.13%'
    'G14' = '10'
    'D15' = '0.001294'
    'E15' = '-0.84%'
    'G15' = '10'
    'D16' = '0.04222'
    'E16' = '-4.23%'
    'G16' = '10'
    'D17' = '0.005924'
    'E17' = '0.00%'
    'G17' = '10'
    'D18' = '3.601'
    'E18' = '2.86%'
    'G18' = '10'
    'D19' = '0.3496'
    'E19' = '-0.25%'
    'G19' = '10'
    'D20' = '8.433'
    'E20' = '-3.71%'
    'G20' = '10'
    'E21' = '0.55%'
    'G21' = '10'
    'G22' = '10'
    'D23' = '0.001246'
    'E23' = '-1.42%'
    'G23' = '10'
    'D24' = '0.004461'
    'E24' = '2.07%'
    'G24' = '10'
    'D25' = '0.0001203'
    'E25' = '-4.61%'
    'G25' = '10'
    'D26' = '0.0003995'
    'E26' = '-0.05%'
    'G26' = '10'
    'G27' = '10'
    'G28' = '10'
    'G29' = '10'
    'G30' = '10'
    'G31' = '10'
    'G32' = '10'
    'G33' = '10'
    'G34' = '10'
    'G35' = '10'
    'G36' = '10'
    'G37' = '10'
    'D38' = '0.02631'
    'E38' = '-5.96%'
    'G38' = '10'
    'D39' = '0.05471'
    'E39' = '-9.95%'
    'G39' = '10'
    'D40' = '0.007700'
    'E40' = '-3.06%'
    'G40' = '10'
    'E41' = '-2.20%'
    'G41' = '10'
    'D42' = '0.007199'
    'E42' = '-19.81%'
    'G42' = '10'
    'D43' = '0.002084'
    'E43' = '-1.78%'
    'G43' = '10'
    'D44' = '0.008854'
    'E44' = '-14.47%'
    'G44' = '10'
    'D45' = '0.00007086'
    'E45' = '1.02%'
    'G45' = '10'
    'E46' = '-0.05%'
    'G46' = '10'
    'D47' = '0.003558'
    'E47' = '11.31%'
    'G47' = '10'
    'E48' = '-0.05%'
    'G48' = '10'
    'E49' = '-0.05%'
    'G49' = '10'
    'E50' = '-0.05%'
    'G50' = '10'
    'G51' = '10'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"   # force text so numeric-looking strings are not reinterpreted
    $cell.Value = $updates[$addr]
}
